$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and 1h volume change (E) columns with latest scraped values.
# Rows 13/14 swap WrappedEther <-> Polkadot ordering; row 51 replaces PancakeSwap with Aave.

$ws.Range("D2").Value = "27.540.78"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "1.845.56"
$ws.Range("E3").Value = "  -1.73%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -1.28%  "

$ws.Range("D5").Value = "'334.05"
$ws.Range("E5").Value = "  -0.21%  "

$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  -1.17%  "

$ws.Range("D7").Value = "'0.4634"
$ws.Range("E7").Value = "  -1.26%  "

$ws.Range("D8").Value = "'0.3847"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("D9").Value = "'45.96"
$ws.Range("E9").Value = "  -1.90%  "

$ws.Range("D10").Value = "'0.07884"
$ws.Range("E10").Value = "  -0.77%  "

$ws.Range("E11").Value = "  -1.03%  "

$ws.Range("D12").Value = "'21.46"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.955"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.844.39"
$ws.Range("E14").Value = "  -2.45%  "

$ws.Range("D15").Value = "'7.127"
$ws.Range("E15").Value = "  +0.47%  "

$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").Value = "'88.44"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "'0.06671"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").Value = "'0.00001034"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").Value = "'17.09"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("E21").Value = "  -1.13%  "

$ws.Range("D22").Value = "27.539.96"
$ws.Range("E22").Value = "  -1.34%  "

$ws.Range("D23").Value = "'5.389"
$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("D24").Value = "'10.89"

$ws.Range("D25").Value = "'2.311"
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("D26").Value = "'158.61"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("E27").Value = "  -2.08%  "

$ws.Range("D28").Value = "'2.114"
$ws.Range("E28").Value = "  +1.88%  "

$ws.Range("D29").Value = "'5.402"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("D30").Value = "'119.82"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").Value = "'0.9761"
$ws.Range("E31").Value = "  +2.04%  "

$ws.Range("D32").Value = "'0.09395"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("D33").Value = "'3.592"
$ws.Range("E33").Value = "  -1.81%  "

$ws.Range("D34").Value = "'5.301"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").Value = "'1.340"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").Value = "'0.06044"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").Value = "'0.02228"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("D38").Value = "'8.285"
$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").Value = "'1.178"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").Value = "'0.5896"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("E41").Value = "  -1.58%  "

$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("D43").Value = "'1.233"
$ws.Range("E43").Value = "  -3.01%  "

$ws.Range("D44").Value = "'0.5584"
$ws.Range("E44").Value = "  -1.06%  "

$ws.Range("D45").Value = "'12.15"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").Value = "'1.909"
$ws.Range("E46").Value = "  -0.47%  "

$ws.Range("D47").Value = "'0.06697"
$ws.Range("E47").Value = "  -2.32%  "

$ws.Range("D48").Value = "'110.90"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("D49").Value = "'1.050"
$ws.Range("E49").Value = "  -1.33%  "

$ws.Range("D50").Value = "'1.005"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'69.95"
$ws.Range("E51").Value = "  -0.88%  "

